# feat: add 2022-Q1 data
#
#  - the current "总计" sheet is repurposed in place into the new "2022-Q1"
#    per-fund holdings sheet (it keeps its original sheetId/r:id identity)
#  - a duplicate of it, made while it still holds the "总计" summary
#    content/style, becomes the refreshed "总计" sheet: the pre-existing
#    2020-Q4 row is pushed down to row 3 and a new 2022-Q1 row is written
#    at row 2

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(2)

# Duplicate "总计" while it is still the summary sheet; the copy is inserted
# immediately after it and becomes the refreshed "总计" sheet.
$wsTotal.Copy($null, $wsTotal)
$wsNewTotal = $wb.ActiveSheet

# The original sheet object is repurposed into the new "2022-Q1" fund sheet.
$wsQ1 = $wsTotal
$wsQ1.Name = "2022-Q1"
$wsNewTotal.Name = "总计"

# ---- "2022-Q1" sheet: per-fund holdings (same layout as "2020-Q4") --------
# Extend the header style (currently only on B1:D1) across the new E1:H1
# columns before writing their text.
$wsQ1.Range("D1").Copy($wsQ1.Range("E1:H1"))

$wsQ1.Cells.Item(1,2).Value = "基金代码"
$wsQ1.Cells.Item(1,3).Value = "基金名称"
$wsQ1.Cells.Item(1,4).Value = "基金规模"
$wsQ1.Cells.Item(1,5).Value = "股票总仓位"
$wsQ1.Cells.Item(1,6).Value = "仓位占比"
$wsQ1.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ1.Cells.Item(1,8).Value = "仓位排名"

$wsQ1.Cells.Item(2,1).Value = 0
$wsQ1.Cells.Item(2,2).Value = "'519779"
$wsQ1.Cells.Item(2,2).ClearFormats()
$wsQ1.Cells.Item(2,3).Value = "交银施罗德沪港深价值精选灵活配置混合"
$wsQ1.Cells.Item(2,4).Value = "'5.13"
$wsQ1.Cells.Item(2,4).ClearFormats()
$wsQ1.Cells.Item(2,5).Value = "'84.44"
$wsQ1.Cells.Item(2,5).ClearFormats()
$wsQ1.Cells.Item(2,6).Value = "'4.64"
$wsQ1.Cells.Item(2,6).ClearFormats()
$wsQ1.Cells.Item(2,7).Value = "'0.2380"
$wsQ1.Cells.Item(2,7).ClearFormats()
$wsQ1.Cells.Item(2,8).Value = 6

# ---- "总计" sheet: push the existing 2020-Q4 row to row 3, add 2022-Q1 row 2
$wsNewTotal.Cells.Item(2,1).Copy($wsNewTotal.Cells.Item(3,1))
$wsNewTotal.Cells.Item(3,1).Value = 1
$wsNewTotal.Cells.Item(3,2).Value = "2020-Q4"
$wsNewTotal.Cells.Item(3,3).Value = 1
$wsNewTotal.Cells.Item(3,4).Value = 0.24

$wsNewTotal.Cells.Item(2,1).Value = 0
$wsNewTotal.Cells.Item(2,2).Value = "2022-Q1"
$wsNewTotal.Cells.Item(2,3).Value = 1
$wsNewTotal.Cells.Item(2,4).Value = 0.24

# Restore the original active sheet/selection.
$wb.Worksheets.Item(1).Activate()
